# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 369
    $ws.Range("F4").Value = 1573
    $ws.Range("F5").Value = 10
    $ws.Range("F7").Value = 403
    $ws.Range("F9").Value = 61
    $ws.Range("F10").Value = 435
}
